{"js": "// Office.js (Word JavaScript API) script\n// Change the list item \"Principais influenciadores\" to\n// \"Principais Influenciadores e Principais Segmentos\" (capitalize the\n// \"I\" in \"Influenciadores\" and append \" e Principais Segmentos\"),\n// reproducing the exact run-split the author's edit left behind:\n//   run1: \"Principais \"\n//   run2: \"I\"\n//   run3: \"nfluenciadores\"\n//   run4: \" e Principais Segmentos\"\n// All four runs keep the original run formatting (rFonts/kern/sz/szCs/\n// lang/ligatures) untouched.\n\nconst body = context.document.body;\n\n// Locate the paragraph holding the exact (case-sensitive) text so we\n// don't collide with the similarly worded sentence earlier in the doc\n// (\"... Narrativa Inteligente, Principais Influenciadores, ...\").\nconst results = body.search(\"Principais influenciadores\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n\n  const rpr =\n    '<w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/>' +\n    '<w:kern w:val=\"0\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:lang w:eastAsia=\"pt-BR\"/>' +\n    '<w14:ligatures w14:val=\"none\"/></w:rPr>';\n\n  let runs = \"\";\n  runs += \"<w:r>\" + rpr + '<w:t xml:space=\"preserve\">Principais </w:t></w:r>';\n  runs += \"<w:r>\" + rpr + \"<w:t>I</w:t></w:r>\";\n  runs += \"<w:r>\" + rpr + \"<w:t>nfluenciadores</w:t></w:r>\";\n  runs += \"<w:r>\" + rpr + '<w:t xml:space=\"preserve\"> e Principais Segmentos</w:t></w:r>';\n\n  const xml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n    \"<w:body><w:p>\" + runs + \"</w:p></w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\";\n\n  target.insertOoxml(xml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop script\n# Change the list item \"Principais influenciadores\" to\n# \"Principais Influenciadores e Principais Segmentos\" (capitalize the\n# \"I\" in \"Influenciadores\" and append \" e Principais Segmentos\"),\n# reproducing the exact run-split the author's edit left behind:\n#   run1: \"Principais \"\n#   run2: \"I\"\n#   run3: \"nfluenciadores\"\n#   run4: \" e Principais Segmentos\"\n# All four runs keep the original run formatting (rFonts/kern/sz/szCs/\n# lang/ligatures) untouched.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph holding the exact (case-sensitive) text so we\n# don't collide with the similarly worded sentence earlier in the doc\n# (\"... Narrativa Inteligente, Principais Influenciadores, ...\").\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Principais influenciadores\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$rng.Find.MatchWildcards = $false\n$found = $rng.Find.Execute()\n\nif ($found) {\n    # Re-seat a plain Range over the hit's [start,end) span \u2014 InsertXML\n    # on the live Find-result range object duplicates the original run\n    # instead of replacing it, so rebuild a fresh Range first.\n    $target = $d.Range($rng.Start, $rng.End)\n\n    $rpr = '<w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/><w:kern w:val=\"0\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:lang w:eastAsia=\"pt-BR\"/><w14:ligatures w14:val=\"none\"/></w:rPr>'\n\n    $runs = ''\n    $runs += '<w:r>' + $rpr + '<w:t xml:space=\"preserve\">Principais </w:t></w:r>'\n    $runs += '<w:r>' + $rpr + '<w:t>I</w:t></w:r>'\n    $runs += '<w:r>' + $rpr + '<w:t>nfluenciadores</w:t></w:r>'\n    $runs += '<w:r>' + $rpr + '<w:t xml:space=\"preserve\"> e Principais Segmentos</w:t></w:r>'\n\n    $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n        '<w:body><w:p>' + $runs + '</w:p></w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n\n    $target.InsertXML($xml, \"Replace\")\n}\n"}
